$p = $ppt.ActivePresentation

# ---- Slide 1: merge the two subtitle runs into one ----
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$trSub = $subtitle.TextFrame.TextRange
$trSub.Text = "__tmp__"
$trSub.Text = "I don’t know, you tell me!"

# ---- Slide 2 (new): "What is Quishing?" ----
$s2 = $p.Slides.Add(2, 2)

$title2 = $s2.Shapes.Item(1)
$trTitle2 = $title2.TextFrame.TextRange
$trTitle2.Text = "What is "
$trTitle2.InsertAfter("Quishing")
$trTitle2.InsertAfter("?")

$content2 = $s2.Shapes.Item(2)
$trContent2 = $content2.TextFrame.TextRange
$trContent2.Text = "__tmp__"
$trContent2.Text = "I don’t know, you tell me!"
$trContent2.InsertAfter([char]13)

# ---- Slide 3 (new): "Protection against Phishing" ----
$s3 = $p.Slides.Add(3, 2)

$title3 = $s3.Shapes.Item(1)
$trTitle3 = $title3.TextFrame.TextRange
$trTitle3.Text = "Protection against Phishing"

$content3 = $s3.Shapes.Item(2)
$trContent3 = $content3.TextFrame.TextRange
$trContent3.Text = "__tmp__"
$trContent3.Text = "I don’t know, you tell me!"
$trContent3.InsertAfter([char]13)
